$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: real date value (9/17/2024) + time 9:00 AM + "Yes"
$ws.Range("A2").Value = 45552
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = 0.375
$ws.Range("B2").NumberFormat = "h:mm"
$ws.Range("C2").Value = "Yes"

# Row 3: date stored as text " 9/17/24" (leading space) + time 1:00 PM + "Yes"
$ws.Range("A3").Value = " 9/17/24"
$ws.Range("B3").Value = 0.54166666666666663
$ws.Range("B3").NumberFormat = "h:mm"
$ws.Range("C3").Value = "Yes"

# Row 4: date stored as text " 9/17/24" (leading space) + time 4:00 PM + "Yes"
$ws.Range("A4").Value = " 9/17/24"
$ws.Range("B4").Value = 0.66666666666666663
$ws.Range("B4").NumberFormat = "h:mm"
$ws.Range("C4").Value = "Yes"

$ws.Range("A5").Select() | Out-Null
